# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's results as row 45.
$ws.Range("A45").Value = 45630
$ws.Range("B45").Value = 119
$ws.Range("C45").Value = 99
$ws.Range("D45").Value = 106

# Row 45 is now the last row, so it takes on the "last row" date format
# (date only, no time) that row 44 previously had.
$ws.Range("A45").NumberFormat = "YYYY-MM-DD"

# Row 44 is no longer the last row, so it reverts to the regular date
# format (with time) used by every other data row.
$ws.Range("A44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
